$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 53 - this shifts the existing rows 53..131 down to 54..132
# and extends the used range to A1:R132 (matching the dimension change in the diff).
$ws.Rows(53).Insert()

# Populate the freshly inserted row 53 with the new weekly record.
$ws.Cells.Item(53, 1).Value = 1
$ws.Cells.Item(53, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(53, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(53, 4).Value = 45203
$ws.Cells.Item(53, 5).Value = 15
$ws.Cells.Item(53, 6).Value = 100112040
$ws.Cells.Item(53, 7).Value = "Cilantro"
$ws.Cells.Item(53, 8).Value = "Sin especificar"
$ws.Cells.Item(53, 9).Value = "Primera"
$ws.Cells.Item(53, 10).Value = 300
$ws.Cells.Item(53, 11).Value = 800
$ws.Cells.Item(53, 12).Value = 1000
$ws.Cells.Item(53, 13).Value = 900
$ws.Cells.Item(53, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(53, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(53, 16).Value = 450
$ws.Cells.Item(53, 17).Value = 2
$ws.Cells.Item(53, 18).Value = "Hortaliza"
